$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

function Clear-Cell($row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# Row 25: bug fix removes leftover empty AX/AY/AZ cells
Clear-Cell 25 50
Clear-Cell 25 51
Clear-Cell 25 52

# Row 26
Set-TextCell 26 1 '1'
Set-TextCell 26 2 '1'
Set-TextCell 26 3 '1'
Set-TextCell 26 4 'Branca'
Set-TextCell 26 5 'Masculino'
Set-TextCell 26 6 '1'
Set-TextCell 26 7 '1'
Set-TextCell 26 8 '1'
Set-TextCell 26 9 '1'
Set-TextCell 26 10 '1'
Set-TextCell 26 11 '2000-01-01'
Set-TextCell 26 12 '1'
Set-TextCell 26 13 '1'
Set-TextCell 26 14 '1'
Set-TextCell 26 15 'NÃO'
Set-TextCell 26 16 'NÃO'
Set-TextCell 26 17 'NÃO'
Set-TextCell 26 18 'NÃO'
Set-TextCell 26 19 'NÃO'
Set-TextCell 26 20 'NÃO'
Set-TextCell 26 21 'NÃO'
Set-TextCell 26 22 'NÃO'
Set-TextCell 26 23 'NÃO'
Set-TextCell 26 24 'NÃO'
Set-TextCell 26 25 'NÃO'
Set-TextCell 26 26 'NÃO'
Set-TextCell 26 27 'NÃO'
Set-TextCell 26 28 'NÃO'
Set-TextCell 26 29 '1'
Set-TextCell 26 30 '1'
Set-TextCell 26 31 '1'
Set-TextCell 26 32 '1'
Set-TextCell 26 33 '1'
Set-TextCell 26 34 'Rural'
Set-TextCell 26 35 '1'
Set-TextCell 26 36 '1'
Set-TextCell 26 37 '1'
Set-TextCell 26 38 '1'
Set-TextCell 26 39 '1'
Set-TextCell 26 40 '1'
Set-TextCell 26 41 '1/1/2000'
Set-TextCell 26 42 '1'
Set-TextCell 26 43 '1/1/2000'
Set-TextCell 26 44 'Manhã'
Set-TextCell 26 45 '01. Berçário I'
Set-TextCell 26 46 '01 - Do Lar'
Set-TextCell 26 47 'NÃO'
Set-TextCell 26 48 'NÃO'
Set-TextCell 26 49 '1'
Set-TextCell 26 50 '1'
Set-TextCell 26 51 '1'
Set-TextCell 26 52 '1'
Set-TextCell 26 53 '2000-01-01'

# Row 27
Set-TextCell 27 1 '1'
Set-TextCell 27 2 'aluno 1'
Set-TextCell 27 3 '123'
Set-TextCell 27 4 'Branca'
Set-TextCell 27 5 'Masculino'
Set-TextCell 27 6 'aa'
Set-TextCell 27 7 'a'
Set-TextCell 27 8 'aa'
Set-TextCell 27 9 'a'
Set-TextCell 27 10 'a'
Set-TextCell 27 11 '2000-01-01'
Set-TextCell 27 12 'aa'
Set-TextCell 27 13 'aa'
Set-TextCell 27 14 '11'
Set-TextCell 27 15 'NÃO'
Set-TextCell 27 16 'NÃO'
Set-TextCell 27 17 'NÃO'
Set-TextCell 27 18 'NÃO'
Set-TextCell 27 19 'NÃO'
Set-TextCell 27 20 'NÃO'
Set-TextCell 27 21 'NÃO'
Set-TextCell 27 22 'NÃO'
Set-TextCell 27 23 'NÃO'
Set-TextCell 27 24 'NÃO'
Set-TextCell 27 25 'NÃO'
Set-TextCell 27 26 'NÃO'
Set-TextCell 27 27 'NÃO'
Set-TextCell 27 28 'NÃO'
Set-TextCell 27 29 'e'
Set-TextCell 27 30 'a'
Set-TextCell 27 31 '1'
Set-TextCell 27 32 'a'
Set-TextCell 27 33 '1'
Set-TextCell 27 34 'Urbana'
Set-TextCell 27 35 '1'
Set-TextCell 27 36 'a'
Set-TextCell 27 37 'pai'
Set-TextCell 27 38 'mae'
Set-TextCell 27 39 'a'
Set-TextCell 27 40 '1'
Set-TextCell 27 41 '1/1/2000'
Set-TextCell 27 42 '1'
Set-TextCell 27 43 '19/1/2000'
Set-TextCell 27 44 'Tempo Integral'
Set-TextCell 27 45 '03. Maternal I'
Set-TextCell 27 46 '04 - Escola Particular'
Set-TextCell 27 47 'NÃO'
Set-TextCell 27 48 'NÃO'
Set-TextCell 27 50 '1'
Set-TextCell 27 51 '2'
Set-TextCell 27 52 '3'
Set-TextCell 27 53 '2000-01-01'

# Row 28
Set-TextCell 28 1 '2'
Set-TextCell 28 2 '123'
Set-TextCell 28 3 '213'
Set-TextCell 28 4 'Branca'
Set-TextCell 28 5 'Masculino'
Set-TextCell 28 6 '1'
Set-TextCell 28 7 '1'
Set-TextCell 28 8 '1'
Set-TextCell 28 9 '1'
Set-TextCell 28 10 '1'
Set-TextCell 28 11 '2000-01-01'
Set-TextCell 28 12 '1'
Set-TextCell 28 13 'aa'
Set-TextCell 28 14 '1'
Set-TextCell 28 15 'NÃO'
Set-TextCell 28 16 'NÃO'
Set-TextCell 28 17 'NÃO'
Set-TextCell 28 18 'NÃO'
Set-TextCell 28 19 'NÃO'
Set-TextCell 28 20 'NÃO'
Set-TextCell 28 21 'NÃO'
Set-TextCell 28 22 'NÃO'
Set-TextCell 28 23 'NÃO'
Set-TextCell 28 24 'NÃO'
Set-TextCell 28 25 'NÃO'
Set-TextCell 28 26 'NÃO'
Set-TextCell 28 27 'NÃO'
Set-TextCell 28 28 'NÃO'
Set-TextCell 28 29 '1'
Set-TextCell 28 30 '1'
Set-TextCell 28 31 '1'
Set-TextCell 28 32 '1'
Set-TextCell 28 33 '1'
Set-TextCell 28 34 'Rural'
Set-TextCell 28 35 '1'
Set-TextCell 28 36 '1'
Set-TextCell 28 37 '1'
Set-TextCell 28 38 '1'
Set-TextCell 28 39 '1'
Set-TextCell 28 40 '1'
Set-TextCell 28 41 '1/1/2000'
Set-TextCell 28 42 '1'
Set-TextCell 28 43 '20/1/2000'
Set-TextCell 28 44 'Tempo Integral'
Set-TextCell 28 45 '08. Ciclo I - 2° Ano'
Set-TextCell 28 46 '06 - Escola Comunitária'
Set-TextCell 28 47 'NÃO'
Set-TextCell 28 48 'NÃO'
Set-TextCell 28 49 '1'
Clear-Cell 28 50
Clear-Cell 28 51
Clear-Cell 28 52
Set-TextCell 28 53 '2000-01-01'

Write-Output "edit applied"